$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Default green")
$ws.Range("B2").Value = 11.68831168831169
$ws.Range("B4").Value = 17.30769230769231
$ws.Range("B8").Value = 15.13157894736842
$ws.Range("B9").Value = 20.98765432098765
$ws.Range("B11").Value = 1.807228915662651

$ws = $wb.Worksheets.Item("Green")
$ws.Range("B2").Value = 34.41558441558442
$ws.Range("B3").Value = 46.98795180722892
$ws.Range("B4").Value = 23.07692307692308
$ws.Range("B5").Value = 49.68152866242038
$ws.Range("B6").Value = 29.51807228915663
$ws.Range("B7").Value = 22.52747252747253
$ws.Range("B8").Value = 36.84210526315789
$ws.Range("B9").Value = 8.024691358024691
$ws.Range("B10").Value = 14.19753086419753
$ws.Range("B11").Value = 22.28915662650602

$ws = $wb.Worksheets.Item("Yellow")
$ws.Range("B2").Value = 27.27272727272727
$ws.Range("B3").Value = 36.74698795180723
$ws.Range("B4").Value = 25.64102564102564
$ws.Range("B5").Value = 17.19745222929936
$ws.Range("B6").Value = 24.69879518072289
$ws.Range("B7").Value = 36.26373626373626
$ws.Range("B8").Value = 19.73684210526316
$ws.Range("B9").Value = 20.37037037037037
$ws.Range("B10").Value = 51.85185185185185
$ws.Range("B11").Value = 25.30120481927711

$ws = $wb.Worksheets.Item("Orange")
$ws.Range("B2").Value = 9.740259740259742
$ws.Range("B3").Value = 5.421686746987952
$ws.Range("B4").Value = 12.17948717948718
$ws.Range("B5").Value = 22.29299363057325
$ws.Range("B6").Value = 19.27710843373494
$ws.Range("B7").Value = 16.48351648351648
$ws.Range("B8").Value = 9.868421052631579
$ws.Range("B9").Value = 22.22222222222222
$ws.Range("B10").Value = 22.8395061728395
$ws.Range("B11").Value = 21.68674698795181

$ws = $wb.Worksheets.Item("Brown")
$ws.Range("B2").Value = 2.597402597402598
$ws.Range("B3").Value = 7.228915662650602
$ws.Range("B4").Value = 5.76923076923077
$ws.Range("B5").Value = 6.369426751592357
$ws.Range("B6").Value = 9.036144578313253
$ws.Range("B7").Value = 9.89010989010989
$ws.Range("B8").Value = 3.947368421052631
$ws.Range("B9").Value = 8.641975308641975
$ws.Range("B10").Value = 8.641975308641975
$ws.Range("B11").Value = 12.65060240963855

$ws = $wb.Worksheets.Item("Red")
$ws.Range("B2").Value = 7.792207792207792
$ws.Range("B3").Value = 3.614457831325301
$ws.Range("B4").Value = 8.333333333333332
$ws.Range("B5").Value = 4.458598726114649
$ws.Range("B6").Value = 13.85542168674699
$ws.Range("B7").Value = 12.63736263736264
$ws.Range("B8").Value = 9.210526315789473
$ws.Range("B9").Value = 12.34567901234568
$ws.Range("B10").Value = 2.469135802469136
$ws.Range("B11").Value = 12.65060240963855

$ws = $wb.Worksheets.Item("Default Red")
$ws.Range("B2").Value = 6.493506493506493
$ws.Range("B4").Value = 7.692307692307693
$ws.Range("B6").Value = 3.614457831325301
$ws.Range("B7").Value = 2.197802197802198
$ws.Range("B8").Value = 5.263157894736842
$ws.Range("B9").Value = 7.407407407407407
$ws.Range("B11").Value = 3.614457831325301

$ws = $wb.Worksheets.Item("Blue")
$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 0
$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = 0
$ws.Range("B7").Value = 0
$ws.Range("B8").Value = 0
$ws.Range("B9").Value = 0
$ws.Range("B10").Value = 0
$ws.Range("B11").Value = 0
